$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting (style) from the row above down into the new row,
# so the new row's date cell (G8) keeps the same date number format.
$ws.Range("A7:H7").Copy()
$ws.Range("A8:H8").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Now set the actual values for the new row 8
$ws.Range("A8").Value = 9783.61
$ws.Range("B8").Value = 9918.5
$ws.Range("C8").Value = 19.170000000000002
$ws.Range("D8").Value = 19.43
$ws.Range("E8").Value = $true
$ws.Range("F8").Value = 1.36
$ws.Range("G8").Value = 42609.488738425927
$ws.Range("H8").Value = $false
